$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5721.2856
$ws.Range("I40").Value = 4099.8184
$ws.Range("K40").Value = 4099.8184
$ws.Range("M40").Value = -3924.8184
$ws.Range("H41").Value = 2480.8333
$ws.Range("I41").Value = 4900
$ws.Range("J41").Value = 1271.25
$ws.Range("K41").Value = 4900
$ws.Range("L41").Value = 1271.25
$ws.Range("M41").Value = -4460
$ws.Range("N41").Value = -2151.25
$ws.Range("H43").Value = 4551
$ws.Range("J43").Value = 4345.8335
$ws.Range("L43").Value = 4345.8335
$ws.Range("N43").Value = -4483.8335
$ws.Range("H76").Value = 3862.0908
$ws.Range("I76").Value = 3397.875
$ws.Range("J76").Value = 5100
$ws.Range("K76").Value = 3397.875
$ws.Range("L76").Value = 5100
$ws.Range("M76").Value = -3082.875
$ws.Range("N76").Value = -5730
$ws.Range("H79").Value = 3862.0908
$ws.Range("I79").Value = 3397.875
$ws.Range("J79").Value = 5100
$ws.Range("K79").Value = 3397.875
$ws.Range("L79").Value = 5100
$ws.Range("M79").Value = -2305.875
$ws.Range("N79").Value = -7284
$ws.Range("H134").Value = 104998
$ws.Range("J134").Value = 104998
$ws.Range("L134").Value = 104998
$ws.Range("N134").Value = -115138
$ws.Range("H135").Value = 29412192
$ws.Range("I135").Value = 38462016
$ws.Range("J135").Value = 260
$ws.Range("K135").Value = 346158144
$ws.Range("L135").Value = 2340
$ws.Range("M135").Value = -346155609
$ws.Range("N135").Value = -7410
$ws.Range("H138").Value = 3442.9524
$ws.Range("I138").Value = 2091.0732
$ws.Range("K138").Value = 6273.219599999999
$ws.Range("M138").Value = -1133.219599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1551.3726
$ws.Range("I32").Value = 1350.1305
$ws.Range("K32").Value = 1350.1305
$ws.Range("M32").Value = -1063.1305
$ws.Range("H74").Value = 45458988
$ws.Range("I74").Value = 58827544
$ws.Range("K74").Value = 58827544
$ws.Range("M74").Value = -58826670
$ws.Range("H77").Value = 45458988
$ws.Range("I77").Value = 58827544
$ws.Range("K77").Value = 294137720
$ws.Range("M77").Value = -294133352
$ws.Range("H102").Value = 1499.4445
$ws.Range("J102").Value = 1097.6666
$ws.Range("L102").Value = 1097.6666
$ws.Range("N102").Value = -4341.6666
$ws.Range("H104").Value = 500000
$ws.Range("J104").Value = 500000
$ws.Range("L104").Value = 500000
$ws.Range("N104").Value = -506988
$ws.Range("H110").Value = 200959.6
$ws.Range("I110").Value = 250824.5
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 250824.5
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = -248779.5
$ws.Range("N110").Value = -5590
$ws.Range("H122").Value = 3736.1333
$ws.Range("I122").Value = 2985.1428
$ws.Range("K122").Value = 8955.428400000001
$ws.Range("M122").Value = -6505.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1471.2667
$ws.Range("I20").Value = 1509.625
$ws.Range("J20").Value = 1427.4286
$ws.Range("K20").Value = 1509.625
$ws.Range("L20").Value = 1427.4286
$ws.Range("M20").Value = -1262.625
$ws.Range("N20").Value = -1921.4286
$ws.Range("H62").Value = 81900
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 81900
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H80").Value = 713.0769
$ws.Range("I80").Value = 717.25
$ws.Range("J80").Value = 711.2222
$ws.Range("K80").Value = 717.25
$ws.Range("L80").Value = 711.2222
$ws.Range("M80").Value = 280.75
$ws.Range("N80").Value = -2707.2222
$ws.Range("H83").Value = 713.0769
$ws.Range("I83").Value = 717.25
$ws.Range("J83").Value = 711.2222
$ws.Range("K83").Value = 3586.25
$ws.Range("L83").Value = 3556.111
$ws.Range("M83").Value = 1405.75
$ws.Range("N83").Value = -13540.111
$ws.Range("H134").Value = 27112558
$ws.Range("I134").Value = 27112558
$ws.Range("K134").Value = 81337674
$ws.Range("M134").Value = -81335139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2279.6
$ws.Range("I16").Value = 1599.5
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 1599.5
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -1312.5
$ws.Range("N16").Value = -5574
$ws.Range("H31").Value = 2709.4912
$ws.Range("I31").Value = 1880.1628
$ws.Range("J31").Value = 5256.7144
$ws.Range("K31").Value = 1880.1628
$ws.Range("L31").Value = 5256.7144
$ws.Range("M31").Value = -1585.1628
$ws.Range("N31").Value = -5846.7144
$ws.Range("H34").Value = 2709.4912
$ws.Range("I34").Value = 1880.1628
$ws.Range("J34").Value = 5256.7144
$ws.Range("K34").Value = 1880.1628
$ws.Range("L34").Value = 5256.7144
$ws.Range("M34").Value = -1678.1628
$ws.Range("N34").Value = -5660.7144
$ws.Range("H99").Value = 3388.7896
$ws.Range("I99").Value = 3187.4443
$ws.Range("J99").Value = 3570
$ws.Range("K99").Value = 3187.4443
$ws.Range("L99").Value = 3570
$ws.Range("M99").Value = -1689.4443
$ws.Range("N99").Value = -6566
$ws.Range("H107").Value = 28695
$ws.Range("I107").Value = 672.5862
$ws.Range("J107").Value = 144787.86
$ws.Range("K107").Value = 672.5862
$ws.Range("L107").Value = 144787.86
$ws.Range("M107").Value = 1247.4138
$ws.Range("N107").Value = -148627.86
$ws.Range("H113").Value = 2279.6
$ws.Range("I113").Value = 1599.5
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1599.5
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 570.5
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 3540.25
$ws.Range("I122").Value = 3498.4546
$ws.Range("K122").Value = 10495.3638
$ws.Range("M122").Value = -8045.363799999999
$ws.Range("H126").Value = 3388.7896
$ws.Range("I126").Value = 3187.4443
$ws.Range("J126").Value = 3570
$ws.Range("K126").Value = 9562.332900000001
$ws.Range("L126").Value = 10710
$ws.Range("M126").Value = -7092.332900000001
$ws.Range("N126").Value = -15650
$ws.Range("H132").Value = 25001394
$ws.Range("I132").Value = 29413178
$ws.Range("J132").Value = 1284.8334
$ws.Range("K132").Value = 88239534
$ws.Range("L132").Value = 3854.5002
$ws.Range("M132").Value = -88237004
$ws.Range("N132").Value = -8914.5002
$ws.Range("H134").Value = 11907366
$ws.Range("I134").Value = 16668796
$ws.Range("J134").Value = 3791.1667
$ws.Range("K134").Value = 50006388
$ws.Range("L134").Value = 11373.5001
$ws.Range("M134").Value = -50003853
$ws.Range("N134").Value = -16443.5001
$ws.Range("H138").Value = 103695
$ws.Range("J138").Value = 103695
$ws.Range("L138").Value = 103695
$ws.Range("N138").Value = -113975

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 380.85715
$ws.Range("I23").Value = 223
$ws.Range("J23").Value = 538.7143
$ws.Range("K23").Value = 669
$ws.Range("L23").Value = 1616.1429
$ws.Range("M23").Value = -434
$ws.Range("N23").Value = -2086.1429
$ws.Range("H131").Value = 1518.7931
$ws.Range("I131").Value = 1026.2941
$ws.Range("K131").Value = 3078.8823
$ws.Range("M131").Value = 1961.1177

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1963.619
$ws.Range("I97").Value = 1999.125
$ws.Range("J97").Value = 1850
$ws.Range("K97").Value = 1999.125
$ws.Range("L97").Value = 1850
$ws.Range("M97").Value = -1503.125
$ws.Range("N97").Value = -2842
$ws.Range("H102").Value = 4922.185
$ws.Range("I102").Value = 3262.5833
$ws.Range("J102").Value = 18199
$ws.Range("K102").Value = 3262.5833
$ws.Range("L102").Value = 18199
$ws.Range("M102").Value = -1640.5833
$ws.Range("N102").Value = -21443
$ws.Range("H122").Value = 6910.087
$ws.Range("I122").Value = 4733.5
$ws.Range("J122").Value = 11885.143
$ws.Range("K122").Value = 14200.5
$ws.Range("L122").Value = 35655.429
$ws.Range("M122").Value = -11750.5
$ws.Range("N122").Value = -40555.429
$ws.Range("H126").Value = 3473.4
$ws.Range("I126").Value = 3594.5
$ws.Range("J126").Value = 2989
$ws.Range("K126").Value = 10783.5
$ws.Range("L126").Value = 8967
$ws.Range("M126").Value = -8313.5
$ws.Range("N126").Value = -13907

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6232
$ws.Range("I68").Value = 6898
$ws.Range("K68").Value = 6898
$ws.Range("M68").Value = -6149
$ws.Range("H71").Value = 6232
$ws.Range("I71").Value = 6898
$ws.Range("K71").Value = 34490
$ws.Range("M71").Value = -30746
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H130").Value = 59999.75
$ws.Range("J130").Value = 59999.75
$ws.Range("L130").Value = 59999.75
$ws.Range("N130").Value = -70039.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2574.375
$ws.Range("I100").Value = 3619.2
$ws.Range("J100").Value = 833
$ws.Range("K100").Value = 7238.4
$ws.Range("L100").Value = 1666
$ws.Range("M100").Value = -6697.4
$ws.Range("N100").Value = -2748
$ws.Range("H113").Value = 1611.0834
$ws.Range("I113").Value = 1467.6666
$ws.Range("J113").Value = 2041.3334
$ws.Range("K113").Value = 4402.9998
$ws.Range("L113").Value = 6124.0002
$ws.Range("M113").Value = -2232.9998
$ws.Range("N113").Value = -10464.0002
$ws.Range("H126").Value = 1160.5834
$ws.Range("I126").Value = 1003
$ws.Range("K126").Value = 3009
$ws.Range("M126").Value = -539
$ws.Range("H135").Value = 91715
$ws.Range("J135").Value = 91715
$ws.Range("L135").Value = 91715
$ws.Range("N135").Value = -101855
